# Add a new "Select" / "Auswählen" translation row to both the "en" and
# "de" resource sheets (row 23), matching the new shared strings that a
# "Select" dropdown placeholder needs (commit: "New scss from resumee").

$wb = $excel.ActiveWorkbook

# --- "en" sheet (Name/Value resource list) -----------------------------
$wsEn = $wb.Worksheets.Item("en")

# Duplicate the formatting of the last existing row (style s="1") onto the
# new row 23, then overwrite its values with the new resource strings.
$wsEn.Range("A22:B22").Copy($wsEn.Range("A23:B23"))
$wsEn.Range("A23").Value = "Select"
$wsEn.Range("B23").Value = "Select"

# --- "de" sheet (German translations) -----------------------------------
$wsDe = $wb.Worksheets.Item("de")

$wsDe.Range("A22:B22").Copy($wsDe.Range("A23:B23"))
$wsDe.Range("A23").Value = "Select"
$wsDe.Range("B23").Value = "Auswählen"

# Leave the "en" sheet active/selected, with the cursor resting below the
# data (matching the saved view state of the source workbook).
$wsEn.Range("A29").Select()
